$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add the new row of data (row 5) ---
$ws.Range("A5").Value = 2019
$ws.Range("B5").Formula = "=B4+54"
$ws.Range("C5").Formula = "=C4+0"
$ws.Range("D5").Formula = "=D4+0"

# --- 2. Column widths ---
$ws.Columns.Item(2).ColumnWidth = 17.42578125
$ws.Columns.Item(3).ColumnWidth = 19.42578125
$ws.Columns.Item(4).ColumnWidth = 16.140625

# --- 3. Convert range into an Excel Table ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:D5"), 0, 1)
$tbl.Name = "Tableau2"
$tbl.TableStyle = "TableStyleLight2"

# --- 4. Selection ---
$ws.Range("A1:D5").Select()

Write-Host "done"
